$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" row (original row 26). All rows below shift up by one.
$ws.Rows.Item(26).Delete()

# Remove the "SC 92" row. After the first deletion it now sits at row 27.
$ws.Rows.Item(27).Delete()

# "SC 119" (now at row 29) loses its previously-imputed value in column B/C.
$ws.Range("C29").ClearContents()

# "SC 232" (now at row 33) gains a value in column B/C where it was previously missing.
$ws.Range("C33").Value = 10.4
